# Add German ("de") translation column F, mirroring the "source" column E
# (the translations simply duplicate the source text, matching the
# not-yet-translated placeholder pattern already used for rows 2-6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use E3's formatting (style index 3 - the standard translated-cell style)
# as the template for the new F column cells, then overwrite each cell's
# value with the corresponding row's column E value.
$ws.Range("E3").Copy()

$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = $ws.Range("E2").Value()

$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = $ws.Range("E3").Value()

$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Value = $ws.Range("E4").Value()

$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value = $ws.Range("E5").Value()

$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F6").Value = $ws.Range("E6").Value()
